$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.020574846168834
$ws.Range("D2").Value = 1.026259465027548
$ws.Range("E2").Value = 1.02155833181374
$ws.Range("F2").Value = 1.031506699814445
$ws.Range("I2").Value = 1.029430903032845
$ws.Range("J2").Value = 1.025771184298331
$ws.Range("K2").Value = 1.029082633784605
$ws.Range("L2").Value = 1.024395307221931
$ws.Range("M2").Value = 1.034314619661036
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.021560697842816
$ws.Range("D3").Value = 1.026984683481957
$ws.Range("E3").Value = 1.02239584836003
$ws.Range("F3").Value = 1.03270218394147
$ws.Range("I3").Value = 1.029600306709796
$ws.Range("J3").Value = 1.026394092810482
$ws.Range("K3").Value = 1.029615551789267
$ws.Range("L3").Value = 1.025039217348184
$ws.Range("M3").Value = 1.035317645065252
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.022198891291063
$ws.Range("D4").Value = 1.027453956667277
$ws.Range("E4").Value = 1.022938399991488
$ws.Range("F4").Value = 1.03347615548905
$ws.Range("I4").Value = 1.029708574408544
$ws.Range("J4").Value = 1.026796870787447
$ws.Range("K4").Value = 1.029959727298783
$ws.Range("L4").Value = 1.02545585119069
$ws.Range("M4").Value = 1.035966519531227
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.022467254880221
$ws.Range("D5").Value = 1.027651240124347
$ws.Range("E5").Value = 1.02316663682961
$ws.Range("F5").Value = 1.033801632493445
$ws.Range("I5").Value = 1.029753767138342
$ws.Range("J5").Value = 1.026966129773703
$ws.Range("K5").Value = 1.030104260478777
$ws.Range("L5").Value = 1.025630998866819
$ws.Range("M5").Value = 1.036239270636496
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.022512318230925
$ws.Range("D6").Value = 1.027684364915302
$ws.Range("E6").Value = 1.023204967460065
$ws.Range("F6").Value = 1.033856287353315
$ws.Range("I6").Value = 1.029761336245072
$ws.Range("J6").Value = 1.026994545047994
$ws.Range("K6").Value = 1.030128518939313
$ws.Range("L6").Value = 1.025660406609689
$ws.Range("M6").Value = 1.036285064658478
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.022202476915362
$ws.Range("D7").Value = 1.027456592775766
$ws.Range("E7").Value = 1.022941449122333
$ws.Range("F7").Value = 1.033480504137764
$ws.Range("I7").Value = 1.02970917954565
$ws.Range("J7").Value = 1.026799132705604
$ws.Range("K7").Value = 1.029961659180674
$ws.Range("L7").Value = 1.025458191543697
$ws.Range("M7").Value = 1.03597016418495
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.020907960875829
$ws.Range("D8").Value = 1.026504553676439
$ws.Range("E8").Value = 1.021841245048125
$ws.Range("F8").Value = 1.031910633920976
$ws.Range("I8").Value = 1.02948843266502
$ws.Range("J8").Value = 1.025981757939122
$ws.Range("K8").Value = 1.029262871781822
$ws.Range("L8").Value = 1.024612923235102
$ws.Range("M8").Value = 1.034653627979122
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.018629028235771
$ws.Range("D9").Value = 1.024827046526447
$ws.Range("E9").Value = 1.019907350405284
$ws.Range("F9").Value = 1.029147469659691
$ws.Range("I9").Value = 1.029089140896511
$ws.Range("J9").Value = 1.024539275869281
$ws.Range("K9").Value = 1.028026507940745
$ws.Range("L9").Value = 1.023123333902996
$ws.Range("M9").Value = 1.032332564516079
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017111215426554
$ws.Range("D10").Value = 1.023708837850876
$ws.Range("E10").Value = 1.018621368710491
$ws.Range("F10").Value = 1.02730744786471
$ws.Range("I10").Value = 1.028816035772145
$ws.Range("J10").Value = 1.0235761964207
$ws.Range("K10").Value = 1.027198933864809
$ws.Range("L10").Value = 1.022130230277164
$ws.Range("M10").Value = 1.030784394130464
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.016454338434647
$ws.Range("D11").Value = 1.023224683492674
$ws.Range("E11").Value = 1.018065312921482
$ws.Range("F11").Value = 1.02651118600082
$ws.Range("I11").Value = 1.028696143035387
$ws.Range("J11").Value = 1.023158840244995
$ws.Range("K11").Value = 1.026839802498169
$ws.Range("L11").Value = 1.021700202067715
$ws.Range("M11").Value = 1.030113826331613
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016210397172966
$ws.Range("D12").Value = 1.023044853510665
$ws.Range("E12").Value = 1.017858887415823
$ws.Range("F12").Value = 1.026215490224437
$ws.Range("I12").Value = 1.028651363952704
$ws.Range("J12").Value = 1.023003765438975
$ws.Range("K12").Value = 1.02670628753073
$ws.Range("L12").Value = 1.021540469790147
$ws.Range("M12").Value = 1.029864717084059
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016262721042019
$ws.Range("D13").Value = 1.023083427343721
$ws.Range("E13").Value = 1.017903161018545
$ws.Range("F13").Value = 1.026278914756068
$ws.Range("I13").Value = 1.028660980326694
$ws.Range("J13").Value = 1.023037031773723
$ws.Range("K13").Value = 1.026734932256741
$ws.Range("L13").Value = 1.021574732927203
$ws.Range("M13").Value = 1.029918153251786
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016434173110304
$ws.Range("D14").Value = 1.023209818552466
$ws.Range("E14").Value = 1.018048247305572
$ws.Range("F14").Value = 1.026486742244779
$ws.Range("I14").Value = 1.028692446590673
$ws.Range("J14").Value = 1.023146022733652
$ws.Range("K14").Value = 1.026828768504885
$ws.Range("L14").Value = 1.021686998558172
$ws.Range("M14").Value = 1.030093235508729
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016539817295755
$ws.Range("D15").Value = 1.02328769323592
$ws.Range("E15").Value = 1.018137655464019
$ws.Range("F15").Value = 1.026614801050094
$ws.Range("I15").Value = 1.028711801460543
$ws.Range("J15").Value = 1.023213169029312
$ws.Range("K15").Value = 1.026886568549292
$ws.Range("L15").Value = 1.021756169047472
$ws.Range("M15").Value = 1.03020110540815
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017154817488713
$ws.Range("D16").Value = 1.023740970444469
$ws.Range("E16").Value = 1.018658288859746
$ws.Range("F16").Value = 1.027360303186138
$ws.Range("I16").Value = 1.028823958203941
$ws.Range("J16").Value = 1.023603887922038
$ws.Range("K16").Value = 1.027222751710706
$ws.Range("L16").Value = 1.022158769713194
$ws.Range("M16").Value = 1.030828893255118
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.01754068357344
$ws.Range("D17").Value = 1.024025309896808
$ws.Range("E17").Value = 1.018985078548248
$ws.Range("F17").Value = 1.02782806489166
$ws.Range("I17").Value = 1.028893873242118
$ws.Range("J17").Value = 1.023848885776898
$ws.Range("K17").Value = 1.027433420391468
$ws.Range("L17").Value = 1.022411308754402
$ws.Range("M17").Value = 1.0312226341444
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.01776578610511
$ws.Range("D18").Value = 1.024191163835686
$ws.Range("E18").Value = 1.01917576491919
$ws.Range("F18").Value = 1.028100948715533
$ws.Range("I18").Value = 1.028934495592433
$ws.Range("J18").Value = 1.023991756393982
$ws.Range("K18").Value = 1.02755622396241
$ws.Range("L18").Value = 1.022558609718706
$ws.Range("M18").Value = 1.031452277298018
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.017842545908158
$ws.Range("D19").Value = 1.024247716317362
$ws.Range("E19").Value = 1.019240796822228
$ws.Range("F19").Value = 1.028194002929248
$ws.Range("I19").Value = 1.028948319971845
$ws.Range("J19").Value = 1.024040466059724
$ws.Range("K19").Value = 1.027598083914806
$ws.Range("L19").Value = 1.022608835391934
$ws.Range("M19").Value = 1.031530576429512
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.017499280341465
$ws.Range("D20").Value = 1.023994802591184
$ws.Range("E20").Value = 1.018950009295077
$ws.Range("F20").Value = 1.02777787370003
$ws.Range("I20").Value = 1.028886388351433
$ws.Range("J20").Value = 1.023822603182974
$ws.Range("K20").Value = 1.027410825483092
$ws.Range("L20").Value = 1.02238421378278
$ws.Range("M20").Value = 1.031180391467363
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016383683310295
$ws.Range("D21").Value = 1.023172599305058
$ws.Range("E21").Value = 1.018005519736375
$ws.Range("F21").Value = 1.026425540266078
$ws.Range("I21").Value = 1.028683187337295
$ws.Range("J21").Value = 1.023113928993338
$ws.Range("K21").Value = 1.026801139310195
$ws.Range("L21").Value = 1.021653939141898
$ws.Range("M21").Value = 1.030041678994027
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.015682564679571
$ws.Range("D22").Value = 1.022655685531479
$ws.Range("E22").Value = 1.017412366912827
$ws.Range("F22").Value = 1.025575686912428
$ws.Range("I22").Value = 1.028554006067373
$ws.Range("J22").Value = 1.022668067510753
$ws.Range("K22").Value = 1.026417125361156
$ws.Range("L22").Value = 1.021194782964416
$ws.Range("M22").Value = 1.029325548503907
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016054212310899
$ws.Range("D23").Value = 1.022929707399184
$ws.Range("E23").Value = 1.017726743236917
$ws.Range("F23").Value = 1.02602617141469
$ws.Range("I23").Value = 1.028622622080673
$ws.Range("J23").Value = 1.022904454462378
$ws.Range("K23").Value = 1.026620762699266
$ws.Range("L23").Value = 1.02143819050369
$ws.Range("M23").Value = 1.02970519974488
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017517988576589
$ws.Range("D24").Value = 1.024008587518508
$ws.Range("E24").Value = 1.018965855347869
$ws.Range("F24").Value = 1.027800552795603
$ws.Range("I24").Value = 1.028889770939256
$ws.Range("J24").Value = 1.023834479256196
$ws.Range("K24").Value = 1.027421035384137
$ws.Range("L24").Value = 1.022396456836588
$ws.Range("M24").Value = 1.03119947917402
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01921792808667
$ws.Range("D25").Value = 1.025260703701645
$ws.Range("E25").Value = 1.02040673384436
$ws.Range("F25").Value = 1.029861443796837
$ws.Range("I25").Value = 1.029193586869212
$ws.Range("J25").Value = 1.024912445355594
$ws.Range("K25").Value = 1.028346727587479
$ws.Range("L25").Value = 1.023508438825193
$ws.Range("M25").Value = 1.035317645065252
